# Insert a new weekly data record at row 211, pushing the existing
# rows 211..319 down to 212..320 (dimension grows from R319 to R320).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 211..319 down by one row.
$ws.Rows(211).Insert()

# Populate the newly inserted row with the new weekly record.
$row = 211

$ws.Cells.Item($row, 1).Value2  = 11
$ws.Cells.Item($row, 2).Value   = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value   = "Bíobío"
$ws.Cells.Item($row, 4).Value2  = 44873
$ws.Cells.Item($row, 5).Value2  = 8
$ws.Cells.Item($row, 6).Value2  = 100112009
$ws.Cells.Item($row, 7).Value   = "Acelga"
$ws.Cells.Item($row, 8).Value   = "Sin especificar"
$ws.Cells.Item($row, 9).Value   = "Primera"
$ws.Cells.Item($row, 10).Value2 = 250
$ws.Cells.Item($row, 11).Value2 = 600
$ws.Cells.Item($row, 12).Value2 = 650
$ws.Cells.Item($row, 13).Value2 = 630
$ws.Cells.Item($row, 14).Value  = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item($row, 15).Value  = "Región de Ñuble"
$ws.Cells.Item($row, 16).Value2 = 630
$ws.Cells.Item($row, 17).Value2 = 1
$ws.Cells.Item($row, 18).Value  = "Hortaliza"
